$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$value = "Merhaba Dünya"

$ws.Range("A2").Value = $value
$ws.Range("A3").Value = $value
